# Adds the two remaining "factory use" rows (electricity turbine use, heat
# turbine use) to both the "emissions" and "removals" sheets, and
# re-capitalizes the existing "o2 factory use" label to "O2 factory use" on
# both sheets (so all three cells collapse onto a single shared string).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "emissions"
$ws2 = $wb.Worksheets.Item(2)   # "removals"

# ---------------------------------------------------------------------
# 1. Re-capitalize "o2 factory use" -> "O2 factory use" on both sheets
# ---------------------------------------------------------------------
$ws1.Range("A53").Value = "O2 factory use"
$ws2.Range("A52").Value = "O2 factory use"

# ---------------------------------------------------------------------
# 2. "emissions" sheet: new rows 59 & 60
# ---------------------------------------------------------------------

# -- formatting donors: copy number formats from existing cells that
#    already carry the styles the new rows need (General+grey = style of
#    C4; scientific+grey = style of J10) rather than hand-building styles.
$ws1.Range("C4:I4").Copy()
$ws1.Range("C59:I59").PasteSpecial(-4122)
$ws1.Range("C4:I4").Copy()
$ws1.Range("C60:I60").PasteSpecial(-4122)

$ws1.Range("J10").Copy()
$ws1.Range("J59").PasteSpecial(-4122)
$ws1.Range("J10").Copy()
$ws1.Range("J60").PasteSpecial(-4122)

$ws1.Range("L48").Copy()
$ws1.Range("L59").PasteSpecial(-4122)
$ws1.Range("L48").Copy()
$ws1.Range("L60").PasteSpecial(-4122)
$ws1.Range("C4").Copy()
$ws1.Range("L59:L60").PasteSpecial(-4122)
$ws1.Range("C4").Copy()
$ws1.Range("F59:F60").PasteSpecial(-4122)

# -- row 59: electricity turbine use
$ws1.Range("A59").Value = "electricity turbine use"
$ws1.Range("C59").Value = [double]"1.2400000000000001E-4"
$ws1.Range("D59").Value = [double]"5.5600000000000001E-6"
$ws1.Range("E59").Value = [double]"1.2200000000000001E-7"
$ws1.Range("F59").Formula = "=SUM(C59:E59)"
$ws1.Range("H59").Value = [double]"4.26E-4"
$ws1.Range("I59").Value = [double]"1.67E-7"
$ws1.Range("J59").Value = [double]"1.6999999999999999E-9"
$ws1.Range("L59").Formula = "=SUM(H59:J59)"

# -- row 60: heat turbine use
$ws1.Range("A60").Value = "heat turbine use"
$ws1.Range("C60").Value = [double]"2.12E-4"
$ws1.Range("D60").Value = [double]"1.04E-5"
$ws1.Range("E60").Value = [double]"1.8900000000000001E-7"
$ws1.Range("F60").Formula = "=SUM(C60:E60)"
$ws1.Range("H60").Value = [double]"9.6599999999999995E-4"
$ws1.Range("I60").Value = [double]"1.02E-8"
$ws1.Range("J60").Value = [double]"1.25E-9"
$ws1.Range("L60").Formula = "=SUM(H60:J60)"

# -- note column, added after both rows exist so the shared-string table
#    picks up "heat turbine use" (row 60) before this repeated note text,
#    matching the order new strings were appended in the authored edit.
$ws1.Range("O59").Value = "includes methane emissions from nat gas use"
$ws1.Range("O60").Value = "includes methane emissions from nat gas use"

# -- bring the selection to the newly-added rows, matching the authored edit
$ws1.Range("A59:A61").Select()

# ---------------------------------------------------------------------
# 3. "removals" sheet: new rows 58 & 59
# ---------------------------------------------------------------------
$ws2.Range("C50").Copy()
$ws2.Range("C58:C59").PasteSpecial(-4122)

$ws2.Range("A58").Value = "electricity turbine use"
$ws2.Range("C58").Value = [double]"5.5799999999999999E-6"

$ws2.Range("A59").Value = "heat turbine use"
$ws2.Range("C59").Value = [double]"9.7100000000000002E-5"

$ws2.Range("C60").Select()
